$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 11559.8
$ws.Range("I86").Value = 2458.8
$ws.Range("K86").Value = 2458.8
$ws.Range("M86").Value = -1335.8

$ws.Range("H89").Value = 11559.8
$ws.Range("I89").Value = 2458.8
$ws.Range("K89").Value = 12294
$ws.Range("M89").Value = -6678

$ws.Range("H107").Value = 905.36365
$ws.Range("I107").Value = 627.53845
$ws.Range("J107").Value = 1937.2858
$ws.Range("K107").Value = 627.53845
$ws.Range("L107").Value = 1937.2858
$ws.Range("M107").Value = 1292.46155
$ws.Range("N107").Value = -5777.2858

$ws.Range("H132").Value = 33337124
$ws.Range("I132").Value = 38465240
$ws.Range("K132").Value = 115395720
$ws.Range("M132").Value = -115393190

$ws.Range("H138").Value = 3143.9412
$ws.Range("I138").Value = 1824.25
$ws.Range("J138").Value = 3550
$ws.Range("K138").Value = 5472.75
$ws.Range("L138").Value = 10650
$ws.Range("M138").Value = -332.75
$ws.Range("N138").Value = -20930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6565.9414
$ws.Range("I32").Value = 6565.9414
$ws.Range("K32").Value = 6565.9414
$ws.Range("M32").Value = -6278.9414

$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976

$ws.Range("H63").Value = 1655
$ws.Range("I63").Value = 1286.1
$ws.Range("J63").Value = 3499.5
$ws.Range("K63").Value = 1286.1
$ws.Range("L63").Value = 3499.5
$ws.Range("M63").Value = -600.0999999999999
$ws.Range("N63").Value = -4871.5

$ws.Range("H66").Value = 1655
$ws.Range("I66").Value = 1286.1
$ws.Range("J66").Value = 3499.5
$ws.Range("K66").Value = 6430.5
$ws.Range("L66").Value = 17497.5
$ws.Range("M66").Value = -2998.5
$ws.Range("N66").Value = -24361.5

$ws.Range("H102").Value = 499
$ws.Range("I102").Value = 499
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 499
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 1123
$ws.Range("N102").ClearContents()

$ws.Range("H110").Value = 2881
$ws.Range("I110").Value = 982.25
$ws.Range("K110").Value = 982.25
$ws.Range("M110").Value = 1062.75

$ws.Range("H122").Value = 3527.6667
$ws.Range("I122").Value = 3493.9285
$ws.Range("K122").Value = 10481.7855
$ws.Range("M122").Value = -8031.7855

$ws.Range("H132").Value = 5329.4546
$ws.Range("I132").Value = 4762.4
$ws.Range("K132").Value = 14287.2
$ws.Range("M132").Value = -11757.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372

$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864

$ws.Range("H109").Value = 79990
$ws.Range("J109").Value = 79990
$ws.Range("L109").Value = 79990
$ws.Range("N109").Value = -82764

$ws.Range("H129").Value = 71666.664
$ws.Range("I129").Value = 80000
$ws.Range("J129").Value = 55000
$ws.Range("K129").Value = 80000
$ws.Range("L129").Value = 55000
$ws.Range("M129").Value = -75000
$ws.Range("N129").Value = -65000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 133.28572
$ws.Range("J7").Value = 267.2
$ws.Range("L7").Value = 267.2
$ws.Range("N7").Value = -493.2

$ws.Range("H16").Value = 8975
$ws.Range("I16").Value = 2950
$ws.Range("J16").Value = 15000
$ws.Range("K16").Value = 2950
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = -2663
$ws.Range("N16").Value = -15574

$ws.Range("H22").Value = 677
$ws.Range("I22").Value = 525
$ws.Range("K22").Value = 525
$ws.Range("M22").Value = -175

$ws.Range("H113").Value = 8975
$ws.Range("I113").Value = 2950
$ws.Range("J113").Value = 15000
$ws.Range("K113").Value = 2950
$ws.Range("L113").Value = 15000
$ws.Range("M113").Value = -780
$ws.Range("N113").Value = -19340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2999
$ws.Range("J80").Value = 2999
$ws.Range("L80").Value = 8997
$ws.Range("N80").Value = -10869

$ws.Range("H83").Value = 2999
$ws.Range("J83").Value = 2999
$ws.Range("L83").Value = 26991
$ws.Range("N83").Value = -36351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8999.4
$ws.Range("I70").Value = 7999
$ws.Range("J70").Value = 9249.5
$ws.Range("K70").Value = 7999
$ws.Range("L70").Value = 9249.5
$ws.Range("M70").Value = -7729
$ws.Range("N70").Value = -9789.5

$ws.Range("H73").Value = 8999.4
$ws.Range("I73").Value = 7999
$ws.Range("J73").Value = 9249.5
$ws.Range("K73").Value = 7999
$ws.Range("L73").Value = 9249.5
$ws.Range("M73").Value = -7063
$ws.Range("N73").Value = -11121.5

$ws.Range("H102").Value = 2868.4
$ws.Range("I102").Value = 2868.4
$ws.Range("K102").Value = 2868.4
$ws.Range("M102").Value = -1246.4

$ws.Range("H107").Value = 179.625
$ws.Range("I107").Value = 135.75
$ws.Range("J107").Value = 223.5
$ws.Range("K107").Value = 135.75
$ws.Range("L107").Value = 223.5
$ws.Range("M107").Value = 1784.25
$ws.Range("N107").Value = -4063.5

$ws.Range("H122").Value = 50853.5
$ws.Range("I122").Value = 54222.637
$ws.Range("K122").Value = 162667.911
$ws.Range("M122").Value = -160217.911

$ws.Range("H126").Value = 3666.3333
$ws.Range("I126").Value = 3666.3333
$ws.Range("K126").Value = 10998.9999
$ws.Range("M126").Value = -8528.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7699.6665
$ws.Range("I40").Value = 7699.6665
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7699.6665
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7563.6665
$ws.Range("N40").ClearContents()

$ws.Range("H93").Value = 884.6667
$ws.Range("I93").Value = 827
$ws.Range("K93").Value = 827
$ws.Range("M93").Value = 421

$ws.Range("H132").Value = 10867.5625
$ws.Range("I132").Value = 13365.25
$ws.Range("K132").Value = 40095.75
$ws.Range("M132").Value = -37565.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2639.2
$ws.Range("I126").Value = 2639.2
$ws.Range("K126").Value = 7917.599999999999
$ws.Range("M126").Value = -5447.599999999999
